$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# Update the repayment strategy value on the input sheet (B17) and restyle it
# (left/top aligned, default font) to match the new "Penalties, Fees, Interest,
# Principal order" option.
$cell = $wsInput.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160

# Move the active tab / selection from the output sheet back to the input
# sheet, with the cursor parked on the cell that was just edited.
$wsInput.Activate()
$wsInput.Range("B17").Select()
